$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 5 (Geetuu Baqalaa) first, then row 2 (Barsiisaa Margsaa),
# so row indices don't shift for the earlier deletion.
$ws.Rows(5).Delete()
$ws.Rows(2).Delete()

# Restore selection as recorded in the post-edit workbook (row 4 selected)
$ws.Range("A4:XFD4").Select()
